# removed button issues with buildpantry and added timetocook filter
#
# Adds the new "Rasgullas" recipe as row 52 of the recipe table on
# Sheet1. Reuses the existing "Indian" cuisine, "green" mark and
# "Desserts" mealtype values already present for the other Indian
# dessert rows, and supplies the new title/description/instructions/
# images/time-to-cook text for this recipe (no "dishtype" value, same
# as the other Indian dessert rows 48-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = "Rasgullas"
$ws.Range("B52").Value = '"If you are intimidated by the idea of making rasgullas at home, please do not be any longer. I was too, until one day I tried it and realized I was worried for no reason all these years. They are quite easy and super quick to make. Give these a try."'
$ws.Range("C52").Value = "`"Prep15 m
Cook35 m
Ready In1 h 50 m
Bring the milk to a boil in a heavy-bottomed pan till it starts foaming; immediately add the lime juice and stir. It will curdle right away. You should see the milk solids (chenna) separate from the whey. Pour into a colander lined with cheesecloth; rinse the chenna with cold water to get rid of the lime juice. Allow the water to drain completely.
Gather the muslin cloth edges like a parcel and express as much water as possible; what you now have is soft paneer. Turn the paneer onto a rolling mat or other smooth surface. Knead the paneer well to make a smooth paste. Roll into a ball and divide into 20 equal portions.
Bring the water to a boil in a pressure cooker; stir the sugar into the boiling water until dissolved.
Roll each portion of paneer into a smooth ball between your palms, making sure there are no cracks; gently drop the balls into the hot syrup. Secure the lid onto the pressure cooker and bring to pressure. Reduce heat to medium-low and pressure cook for 6 minutes.
Release the pressure from the cooker while running under water; remove the lid. The rasgullas should be floating on the syrup and have expanded 2 or 3 times in size. Pour the rasgullas and syrup into a bowl. Gently stir the cardamom into the mixture. Refrigerate to chill completely before serving cold.`""
$ws.Range("D52").Value = "Indian"
$ws.Range("F52").Value = "green"
$ws.Range("G52").Value = '{"https://images.media-allrecipes.com/userphotos/250x250/564413.jpg","https://encrypted-tbn0.gstatic.com/images?q=tbn%3AANd9GcReVDzS5a7OGIQrGskHSvQTVfAJ4ce0iGX8sXUXt2JFLAlSjft7"}'
$ws.Range("H52").Value = "1:50:00"
$ws.Range("I52").Value = "Desserts"

# Setting a multi-line value (C52) made Excel auto-expand the row height
# and pin a custom "ht" on the row; re-running AutoFit restores the
# sheet's default (un-pinned) row height so row 52 matches the styling
# of every other data row in the table.
$ws.Rows.Item(52).AutoFit()
